$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M2").Value = 1.905108
$ws.Range("N2").Value = 3.810216
$ws.Range("O2").Value = 0.07580486173280727
$ws.Range("P2").Value = 0.05715529216076502
$ws.Range("Q2").Value = 1.094556940104
$ws.Range("R2").Value = 6.567341640624
$ws.Range("S2").Value = 0.07580486173280727
$ws.Range("T2").Value = 0.05715529216076502

$ws.Range("O3").Value = 0.1220879833796353
$ws.Range("P3").Value = 0.1380777076800943
$ws.Range("S3").Value = 0.1220879833796353
$ws.Range("T3").Value = 0.1380777076800943

$ws.Range("M4").Value = 3.725954
$ws.Range("N4").Value = 11.177862
$ws.Range("O4").Value = 0.148256911310435
$ws.Range("P4").Value = 0.1676739503331867
$ws.Range("Q4").Value = 2.140702159252
$ws.Range("R4").Value = 19.266319433268
$ws.Range("S4").Value = 0.148256911310435
$ws.Range("T4").Value = 0.1676739503331867

$ws.Range("M5").Value = 6.825836
$ws.Range("N5").Value = 13.651672
$ws.Range("O5").Value = 0.2716022158275637
$ws.Range("P5").Value = 0.2047824327132465
$ws.Range("Q5").Value = 3.921702163768
$ws.Range("R5").Value = 23.530212982608
$ws.Range("S5").Value = 0.2716022158275637
$ws.Range("T5").Value = 0.2047824327132465

$ws.Range("M6").Value = 4.247626666666666
$ws.Range("N6").Value = 12.74288
$ws.Range("O6").Value = 0.1690144349607748
$ws.Range("P6").Value = 0.1911500632430207
$ws.Range("Q6").Value = 2.440422929813333
$ws.Range("R6").Value = 21.96380636832
$ws.Range("S6").Value = 0.1690144349607748
$ws.Range("T6").Value = 0.1911500632430207

$ws.Range("M7").Value = 5.358931000000001
$ws.Range("N7").Value = 16.076793
$ws.Range("O7").Value = 0.213233592788784
$ws.Range("P7").Value = 0.2411605538696867
$ws.Range("Q7").Value = 3.078909498878001
$ws.Range("R7").Value = 27.710185489902
$ws.Range("S7").Value = 0.213233592788784
$ws.Range("T7").Value = 0.2411605538696867
